# fix(publipostage): Add space before ":"
# Also refreshes the underlying clinical-trials data: two rows swap which
# NCT id / title they describe, a new EudraCT id is added, a couple of
# intervention_type values change, and a brand-new trial row is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the colon spacing in "statut_name" text, every row uses the same value ---
$ws.Range("B2").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B3").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B4").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B5").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B6").Value = "4 : pas de résultats postés ni publiés"

# --- Row 3 now describes trial NCT04549194 (was NCT03433859) ---
$ws.Range("C3").Value = "NCT04549194"
$ws.Range("G3").Value = "Contribution of L-Tyrosine to Recovery From Operational Strain on Return From External Operation"
$ws.Range("H3").Value = "USOP"
$ws.Range("I3").Value = "DEVICE"

# --- Row 5 now describes trial NCT03433859 (was NCT04549194), plus a new EudraCT id ---
$ws.Range("C5").Value = "NCT03433859"
$ws.Range("D5").Value = "2014-002068-34"
$ws.Range("G5").Value = "Prospective Multicentric Open Randomised Controlled Trial Comparing Topical Aluminium Chloride to OnabotulinumtoxinA Intradermal Injections in Residual Limb Hyperhidrosis (Lower Limbs)"
$ws.Range("H5").Value = "SALUTOX"
$ws.Range("I5").Value = "OTHER"

# --- Row 6: intervention type updated ---
$ws.Range("I6").Value = "BEHAVIORAL"

# --- New row 7: brand-new trial entry ---
# A7 and F7 hold text that looks numeric ("4" / "2024"); force the Text
# number format first so Excel stores them as text (matching the other
# rows), then restore the default style so no visible formatting changes.
$a7 = $ws.Range("A7")
$f7 = $ws.Range("F7")
$a7.NumberFormat = "@"
$f7.NumberFormat = "@"

$ws.Range("A7").Value = "4"
$ws.Range("B7").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C7").Value = "NCT06174181"
$ws.Range("F7").Value = "2024"
$ws.Range("G7").Value = "Preventive TREatment of Dry Eye in Patients Receiving Repeated Intravitreal Injections for Age-related Macular Degeneration"
$ws.Range("H7").Value = "TREDIA"
$ws.Range("I7").Value = "DEVICE"

$a7.Style = "Normal"
$f7.Style = "Normal"
